# Applies the "intellij shortcuts.docx" edit:
#  1. Remove the _GoBack bookmark that sits after "Shortuts" in the title.
#  2. Rework the "Refractor :- ctrk+shift+L" paragraph into
#     "Refractor  and format:- ctrl+alt+L", dropping the gramStart/gramEnd
#     proof-error markers around "Refractor :" and re-homing the _GoBack
#     bookmark in the middle of the new text ("... and f|ormat:- ...").
#  3. Fix the "ctrk+shift+L" typo -> "ctrl+alt+L" (already folded into the
#     rewritten paragraph XML below).

$d = $word.ActiveDocument

# --- Step 1: drop the stray _GoBack bookmark from the heading paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2 & 3: rebuild the "Refractor" paragraph (4th paragraph) ---
$targetText = "Refractor :- ctrk+shift+L"
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq $targetText) {
        $found = $true

        $rPr = '<w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>'

        $newParagraphXml =
            '<w:p w14:paraId="3B699370" w14:textId="48A98315" w:rsidR="008B7715" w:rsidRDefault="008B7715">' +
                '<w:pPr>' + $rPr + '</w:pPr>' +
                '<w:r>' + $rPr + '<w:t xml:space="preserve">Refractor </w:t></w:r>' +
                '<w:r>' + $rPr + '<w:t xml:space="preserve"> and f</w:t></w:r>' +
                '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
                '<w:bookmarkEnd w:id="0"/>' +
                '<w:r>' + $rPr + '<w:t xml:space="preserve">ormat:- </w:t></w:r>' +
                '<w:proofErr w:type="spellStart"/>' +
                '<w:r>' + $rPr + '<w:t>ctrl+alt+L</w:t></w:r>' +
                '<w:proofErr w:type="spellEnd"/>' +
            '</w:p>'

        $packageXml =
            '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                    '<pkg:xmlData>' +
                        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
                            '<w:body>' + $newParagraphXml + '</w:body>' +
                        '</w:document>' +
                    '</pkg:xmlData>' +
                '</pkg:part>' +
            '</pkg:package>'

        # Replace the whole paragraph (including its end-of-paragraph mark)
        # so no orphan proofErr markers are left behind.
        $para.Range.InsertXML($packageXml)
        break
    }
}

if (-not $found) {
    throw "Could not locate the 'Refractor :- ctrk+shift+L' paragraph"
}

Write-Output "done"
